# Checked on Logical functions and formulas
# Adds a "Meets Both" Goal-Met/Not-Met helper column (P) to the Data sheet,
# unhides the L:N helper columns, and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Unhide columns L:N (12:14) that were previously hidden helper columns.
$ws.Columns.Item(12).Hidden = $false
$ws.Columns.Item(13).Hidden = $false
$ws.Columns.Item(14).Hidden = $false

# L gets a plain fixed width; M and N are widened to fit their header text
# ("Meets Salary" / "Meets Both (1 or 0)").
$ws.Columns.Item(12).ColumnWidth = 8.17
$ws.Columns.Item(13).ColumnWidth = 11
$ws.Columns.Item(14).ColumnWidth = 16.7

# New column P: "Goal Met" / "Not Met" based on column M (Meets Salary)
$ws.Range("P3").Formula = '=IF(M3,"Goal Met","Not Met")'
$ws.Range("P4:P12").Formula = '=IF(M4,"Goal Met","Not Met")'

$ws.Range("S11").Select()
